# Update the "addListItem" sheet: the quick-add test user moves from
# "UserfifteenX" / "ADLILC.8860" to "UserfifteenY" / "ADLILC.8861".
$wb = $excel.ActiveWorkbook

$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "UserfifteenY"
$wsAddListItem.Range("D2").Value = "ADLILC.8861"

# Update the "createUser" sheet: bump the auto-generated user number from
# 1066 to 1068 (dependent formulas in B2/F2 recalc automatically).
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 1068

# The active/selected sheet moves from "createUser" to "addListItem".
$wsAddListItem.Activate()
